$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column T (2022) data, copying the style from column S for each row.
$ws.Range("T4").Value = 2022
$ws.Range("T4").Style = $ws.Range("S4").Style

$ws.Range("T5").Value = 3.7
$ws.Range("T5").Style = $ws.Range("S5").Style

$ws.Range("T6").Value = 1.6
$ws.Range("T6").Style = $ws.Range("S6").Style

$ws.Range("T7").Value = 1.7
$ws.Range("T7").Style = $ws.Range("S7").Style

$ws.Range("T8").Value = 17.9
$ws.Range("T8").Style = $ws.Range("S8").Style

$ws.Range("T9").Value = 7.5
$ws.Range("T9").Style = $ws.Range("S9").Style

$ws.Range("T10").Value = 1.1
$ws.Range("T10").Style = $ws.Range("S10").Style

$ws.Range("T11").Value = 4.4
$ws.Range("T11").Style = $ws.Range("S11").Style

$ws.Range("T12").Value = 3
$ws.Range("T12").Style = $ws.Range("S12").Style

$ws.Range("T13").Value = 4.1
$ws.Range("T13").Style = $ws.Range("S13").Style

$ws.Range("T14").Value = 0.8
$ws.Range("T14").Style = $ws.Range("S14").Style

# Update the selection to match the post-edit state.
$ws.Range("U4").Select()
